$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" (D) and "Data" (E) columns are stored as plain text (inline
# strings) in this workbook, not numbers - e.g. D2 holds the literal text
# "245.69". Excel's default `.Value` setter auto-coerces numeric-looking
# strings into real numbers, which both changes the stored cell type and
# can introduce binary floating-point noise (e.g. "245.82" round-trips as
# 245.81999999999999). To keep the updated cells exactly as text, we
# temporarily force a text number format before assigning the value, then
# restore each cell's original style so no visual/formatting footprint is
# left behind.

$updates = [ordered]@{
    # Price column
    "D2"  = "245.82"
    "D3"  = "23.88"
    "D4"  = "5.207"
    "D5"  = "0.05735"
    "D6"  = "6.479"
    "D7"  = "3.152"
    "D8"  = "0.8144"
    "D9"  = "0.8585"
    "D10" = "0.1378"
    "D11" = "0.06990"
    "D12" = "0.03211"
    "D13" = "0.02877"
    "D14" = "0.09377"
    "D16" = "0.001526"
    "D18" = "0.0005987"
    "D19" = "0.006185"
    "D20" = "0.001243"
    "D21" = "0.004783"
    "D22" = "0.00008497"
    "D23" = "3.528"
    "D24" = "2.153"
    "D25" = "0.3195"
    "D40" = "0.03702"
    "D41" = "0.006387"
    "D44" = "0.007800"
    "D45" = "0.00005488"
    "D47" = "0.3882"
    "D48" = "0.002656"
    "D50" = "0.0001999"

    # Data column
    "E12" = "11LiechtensteinCryptoassetsExchangeLCX"
    "E18" = "17OneONE"
    "E41" = "40KickTokenKICKBestin24h"
    "E48" = "47BOLOBOLOWorstin24h"
}

# Stash each cell's current style and switch to a text number format so
# the assigned value is stored verbatim as a string.
$origStyles = @{}
foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $origStyles[$addr] = $cell.Style
    $cell.NumberFormat = "@"
}

# Assign the new text values, then restore each cell's original style.
foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = $updates[$addr]
    $cell.Style = $origStyles[$addr]
}
